$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Row = 291; Date = 44926; C = 41678378000000; D = 41678378000000; E = 41678378000000; F = 41678378000000; G = 0 },
    @{ Row = 292; Date = 44957; C = 42206471000000; D = 42206471000000; E = 42206471000000; F = 42206471000000; G = 0 },
    @{ Row = 293; Date = 44985; C = 41767936000000; D = 41767936000000; E = 41767936000000; F = 41767936000000; G = 0 },
    @{ Row = 294; Date = 45016; C = 42100727000000; D = 42100727000000; E = 42100727000000; F = 42100727000000; G = 0 },
    @{ Row = 295; Date = 45046; C = 41265420000000; D = 41265420000000; E = 41265420000000; F = 41265420000000; G = 0 },
    @{ Row = 296; Date = 45077; C = 40642084000000; D = 40642084000000; E = 40642084000000; F = 40642084000000; G = 0 },
    @{ Row = 297; Date = 45107; C = 41806284000000; D = 41806284000000; E = 41806284000000; F = 41806284000000; G = 0 }
)

foreach ($r in $newRows) {
    $rowNum = $r.Row

    # Copy the formatting of the last existing data row (290) down onto the
    # new row first, so the date column keeps its custom date/time style.
    $ws.Range("A290:G290").Copy()
    $ws.Range("A$rowNum`:G$rowNum").PasteSpecial(-4122)  # xlPasteFormats

    $ws.Cells.Item($rowNum, 1).Value = $r.Date
    $ws.Cells.Item($rowNum, 2).Value = "ECONOMICS:CNCBBS"
    $ws.Cells.Item($rowNum, 3).Value = $r.C
    $ws.Cells.Item($rowNum, 4).Value = $r.D
    $ws.Cells.Item($rowNum, 5).Value = $r.E
    $ws.Cells.Item($rowNum, 6).Value = $r.F
    $ws.Cells.Item($rowNum, 7).Value = $r.G
}
